# Commit before Merge 5/22
# Applies:
#  - Renames the 4 "lpuser-..." username strings in column E to newly
#    generated test usernames (dated 5-22-2012)
#  - Updates date serials in column B (38144 -> 38165)
#  - Turns C2 (an email address) into a mailto hyperlink with the built-in
#    "Hyperlink" cell style
#  - Updates the saved selection on Sheet1 to P10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the date values in column B (row 2-5): 38144 -> 38165 ---
$ws.Range("B2").Value = 38165
$ws.Range("B3").Value = 38165
$ws.Range("B4").Value = 38165
$ws.Range("B5").Value = 38165

# --- Rename the "lpuser-..." usernames in column E, keeping row mapping ---
$ws.Range("E2").Value = "lpuser-5-22-2012-59800"
$ws.Range("E3").Value = "lpuser-5-22-2012-59840"
$ws.Range("E4").Value = "lpuser-5-22-2012-59876"
$ws.Range("E5").Value = "lpuser-5-22-2012-59911"

# --- Turn C2 into a mailto hyperlink (adds Hyperlink font/style) ---
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:TestFname-a@mailinator.com")

# --- Update the remembered selection for the sheet view ---
$ws.Range("P10").Select()
